$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 24.14949828602258

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 189.6080260415259
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 206.5445078528883

$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 9.983522426115931
$ws.Range("D4").Value = 189.6080260415259
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 214.9010425794622
